# Auto-generated edits applying the "Update countries & provincias Spain" diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 6 de Abril de 2020 a las 08:52'

# Row 4
$ws.Range("B4").Value = 336851
$ws.Range("C4").Value = 178
$ws.Range("E4").Value = 309254
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 9620

# Row 13
$ws.Range("D13").Value = 7298
$ws.Range("E13").Value = 13087

# Row 21
$ws.Range("D21").Value = 585
$ws.Range("E21").Value = 7975

# Row 27
$ws.Range("B27").Value = 4591
$ws.Range("C27").Value = 4
$ws.Range("E27").Value = 4423
$ws.Range("F27").Value = 84
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 72

# Row 31
$ws.Range("D31").Value = 162
$ws.Range("E31").Value = 3846

# Row 43
$ws.Range("A43").Value = 'Finlandia'
$ws.Range("B43").Value = 2176
$ws.Range("C43").Value = 249
$ws.Range("D43").Value = 300
$ws.Range("E43").Value = 1848
$ws.Range("F43").Value = 73
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 28

# Row 44
$ws.Range("A44").Value = 'Mexico'
$ws.Range("B44").Value = 2143
$ws.Range("C44").Value = 253
$ws.Range("D44").Value = 633
$ws.Range("E44").Value = 1416
$ws.Range("F44").Value = 293
$ws.Range("G44").Value = 15
$ws.Range("H44").Value = 94

# Row 45
$ws.Range("A45").Value = 'Panama'
$ws.Range("B45").Value = 1988
$ws.Range("C45").Value = 187
$ws.Range("D45").Value = 13
$ws.Range("E45").Value = 1921
$ws.Range("F45").Value = 78
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 54

# Row 72
$ws.Range("B72").Value = 661
$ws.Range("C72").Value = 7
$ws.Range("D72").Value = 44
$ws.Range("E72").Value = 591
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 26

# Row 82
$ws.Range("D82").Value = 55
$ws.Range("E82").Value = 454

# Row 84
$ws.Range("D84").Value = 8
$ws.Range("E84").Value = 476

# Row 89
$ws.Range("A89").Value = 'Taiwan'
$ws.Range("B89").Value = 373
$ws.Range("C89").Value = 10
$ws.Range("D89").Value = 57
$ws.Range("E89").Value = 311
$ws.Range("H89").Value = 5

# Row 90
$ws.Range("A90").Value = 'Afganistan'
$ws.Range("B90").Value = 367
$ws.Range("C90").Value = 18
$ws.Range("D90").Value = 17
$ws.Range("E90").Value = 343
$ws.Range("H90").Value = 7

# Row 95
$ws.Range("A95").Value = 'Oman'
$ws.Range("B95").Value = 331
$ws.Range("C95").Value = 33
$ws.Range("D95").Value = 61
$ws.Range("E95").Value = 268
$ws.Range("F95").Value = 3
$ws.Range("H95").Value = 2

# Row 96
$ws.Range("A96").Value = 'Cuba'
$ws.Range("B96").Value = 320
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 15
$ws.Range("E96").Value = 297
$ws.Range("F96").Value = 11
$ws.Range("H96").Value = 8

# Row 97
$ws.Range("A97").Value = 'Honduras'
$ws.Range("C97").Value = 30
$ws.Range("D97").Value = 6
$ws.Range("E97").Value = 270
$ws.Range("F97").Value = 10
$ws.Range("H97").Value = 22

# Row 105
$ws.Range("A105").Value = 'Montenegro'
$ws.Range("B105").Value = 223
$ws.Range("C105").Value = 9
$ws.Range("D105").Value = 1
$ws.Range("E105").Value = 220
$ws.Range("F105").Value = 4

# Row 106
$ws.Range("A106").Value = 'Senegal'
$ws.Range("B106").Value = 222
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 82
$ws.Range("E106").Value = 138
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 2

# Row 107
$ws.Range("A107").Value = 'Kirguistan'
$ws.Range("B107").Value = 216
$ws.Range("C107").Value = 69
$ws.Range("D107").Value = 33
$ws.Range("E107").Value = 179
$ws.Range("F107").Value = 5
$ws.Range("G107").Value = 3
$ws.Range("H107").Value = 4

# Row 109
$ws.Range("A109").Value = 'Georgia'
$ws.Range("B109").Value = 188
$ws.Range("C109").Value = 14
$ws.Range("D109").Value = 36
$ws.Range("E109").Value = 150
$ws.Range("F109").Value = 6
$ws.Range("H109").Value = 2

# Row 110
$ws.Range("A110").Value = 'Niger'
$ws.Range("B110").Value = 184
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 13
$ws.Range("E110").Value = 161
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 10

# Row 111
$ws.Range("A111").Value = 'Bolivia'
$ws.Range("B111").Value = 183
$ws.Range("C111").Value = 26
$ws.Range("D111").Value = 2
$ws.Range("E111").Value = 170
$ws.Range("F111").Value = 3
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 11

# Row 112
$ws.Range("A112").Value = 'Islas Feroe'
$ws.Range("B112").Value = 183
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 107
$ws.Range("E112").Value = 76
$ws.Range("F112").Value = 1
$ws.Range("H112").Value = 0

# Row 113
$ws.Range("A113").Value = 'Sri Lanka'
$ws.Range("B113").Value = 176
$ws.Range("D113").Value = 33
$ws.Range("E113").Value = 138
$ws.Range("F113").Value = 5
$ws.Range("H113").Value = 5
